# DaySale_2025-07-12_00-00.xlsx update
# - Adds a new missing-item row ("صابون ديتول اوريجنيال 115 جم") as item #17
#   right before the totals row, pushing the totals row and the footer
#   (timestamp / page / credit) row down by one.
# - Bumps the grand total in column P from 815.4 to 845.4.
# - Refreshes the generation timestamp from 12:12 PM to 12:13 PM.

function Set-TextValue($range, [string]$text) {
    # Force the cell to literal text even when the string looks numeric
    # (e.g. "0", "30.00", "30.0000") so Excel doesn't silently convert it
    # to a Double and drop the formatting/shared-string semantics.
    $origFormat = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.NumberFormat = $origFormat
}

function Copy-BottomBorder($srcRange, $dstRange) {
    # Rows.Insert()/Range.Insert() already clones font/fill/number-format
    # from the row above, but it does not keep the thin bottom-border rule
    # line, so re-apply it explicitly per cell.
    $srcBorder = $srcRange.Borders.Item(9)
    $dstBorder = $dstRange.Borders.Item(9)
    if ($srcBorder.LineStyle -ne $dstBorder.LineStyle) {
        $dstBorder.LineStyle = $srcBorder.LineStyle
        if ($srcBorder.LineStyle -ne -4142) {
            $dstBorder.Weight = $srcBorder.Weight
            $dstBorder.Color = $srcBorder.Color
        }
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

# ---------------------------------------------------------------------
# 1. Insert a new row 23 (shifts old row23 -> 24, old row24 -> 25), using
#    the row above (22) as the formatting template.
# ---------------------------------------------------------------------
$ws.Range("A23:Q23").Insert(-4121)

# ---------------------------------------------------------------------
# 2. Populate the new item row with the new product's data.
# ---------------------------------------------------------------------
$ws.Range("A23").Value2 = 17
Set-TextValue $ws.Range("C23") "صابون ديتول اوريجنيال 115 جم"
Set-TextValue $ws.Range("H23") "9:0"
Set-TextValue $ws.Range("L23") "0"
Set-TextValue $ws.Range("N23") "30.00"
Set-TextValue $ws.Range("P23") "30.0000"
Set-TextValue $ws.Range("Q23") "1:0"

# Re-create the per-cell bottom border lost on insert, and restore the
# row height used by the other item rows.
foreach ($col in $cols) {
    Copy-BottomBorder $ws.Range($col + "22") $ws.Range($col + "23")
}
$ws.Rows.Item(23).RowHeight = 24.75

# Merge the new row's cells the same way the other item rows are merged.
$ws.Range("A23:B23").Merge()
$ws.Range("C23:G23").Merge()
$ws.Range("H23:K23").Merge()
$ws.Range("L23:M23").Merge()
$ws.Range("N23:O23").Merge()

# ---------------------------------------------------------------------
# 3. Totals row (now row 24): bump the grand total and fix its height.
# ---------------------------------------------------------------------
$ws.Range("P24").Value2 = 845.4
$ws.Rows.Item(24).RowHeight = 25.5

# ---------------------------------------------------------------------
# 4. Footer row (now row 25): refresh the generation timestamp.
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("A25") "Saturday, 12 July, 2025 12:13 PM"

Write-Host "Edit applied"
